$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (row 4) — same shape as row 3, new Id 102
$ws.Range("A4").Value = 102
$ws.Range("B4").Value = "note1"
$ws.Range("C4").Value = "note1"
$ws.Range("D4").Value = "NoteTemplate1"
$ws.Range("E4").Value = "[3]"

# Copy style (s="1", centered alignment, default fill) from row 3 to row 4
$ws.Range("A3:E3").Copy()
$ws.Range("A4:E4").PasteSpecial(-4122)  # xlPasteFormats

# Match the new active selection cell recorded in the saved file
$ws.Range("F8").Select() | Out-Null
